# "back to normal code"
# Revert the F11:I12 ratio formulas from "$F$8/<col>8" style back to
# "<col>8/$F$8" style (and analogously for row 9 / row 12), and update the
# active selection on the sheet from I14 to M14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11: ratios against F8 (Max Slow)
$ws.Range("F11").Formula = '=F8/$F$8'
$ws.Range("G11").Formula = '=G8/$F$8'
$ws.Range("H11").Formula = '=H8/$F$8'
$ws.Range("I11").Formula = '=I8/$F$8'

# Row 12: ratios against F9 (Wght Speed)
$ws.Range("F12").Formula = '=F9/$F$9'
$ws.Range("G12").Formula = '=G9/$F$9'
$ws.Range("H12").Formula = '=H9/$F$9'
$ws.Range("I12").Formula = '=I9/$F$9'

# Update the selected cell shown in the sheet view
$ws.Range("M14").Select()
